$wb = $excel.ActiveWorkbook

# --- Update Tire_Type shared string "710Rバフ100" -> "710R" (K2:K6 on Step3_* sheets) ---
$tireTypeSheets = @("Step3_DataPts_0.5","Step3_DataPts_0.7","Step3_DataPts_0.8","Step3_DataPts_0.9")
foreach ($sn in $tireTypeSheets) {
    $s = $wb.Worksheets.Item($sn)
    for ($r = 2; $r -le 6; $r++) {
        $s.Cells.Item($r, 11).Value = "710R"
    }
}

# --- Step1_Data (84 cell updates) ---
$ws = $wb.Worksheets.Item("Step1_Data")
$ws.Range("D2").Value = 0.2286418605659838
$ws.Range("F2").Value = 0.4476256877643685
$ws.Range("H2").Value = 0.02594290329360358
$ws.Range("M2").Value = 0.01534085053127212
$ws.Range("N2").Value = 0.02974716198938816
$ws.Range("O2").Value = 0.07955268351463433
$ws.Range("S2").Value = 0.04696276741955353
$ws.Range("T2").Value = 0.001844424448398699
$ws.Range("V2").Value = 0.008145657074384242
$ws.Range("X2").Value = 0.03004987957805371
$ws.Range("Z2").Value = 0.035713572500001
$ws.Range("AD2").Value = 0.00495476557886092
$ws.Range("AF2").Value = 0.02034017435382522
$ws.Range("AI2").Value = 0.02513761138767244
$ws.Range("D3").Value = 0.2030359476101878
$ws.Range("F3").Value = 0.4671898946942191
$ws.Range("G3").Value = 0.0009979513723979152
$ws.Range("H3").Value = 0.02609948168871957
$ws.Range("M3").Value = 0.01521711615217846
$ws.Range("O3").Value = 0.1049924716440094
$ws.Range("S3").Value = 0.008126883096635232
$ws.Range("T3").Value = 0.02376309779150702
$ws.Range("W3").Value = 0.01217228462583452
$ws.Range("X3").Value = 0.01152589977333775
$ws.Range("Z3").Value = 0.04537013605533151
$ws.Range("AA3").Value = 0.01533219932413948
$ws.Range("AD3").Value = 0.01676121160870259
$ws.Range("AF3").Value = 0.02418365173443643
$ws.Range("AI3").Value = 0.02523177282836344
$ws.Range("D4").Value = 0.1152529212304223
$ws.Range("E4").Value = 0.04413316267050671
$ws.Range("F4").Value = 0.4624364012047548
$ws.Range("G4").Value = 0.05115389306322299
$ws.Range("H4").Value = 0.03923686602650097
$ws.Range("I4").Value = 0.005497491812740805
$ws.Range("O4").Value = 0.08291858815308162
$ws.Range("R4").Value = 0.004449977888918484
$ws.Range("S4").Value = 0.001401171764439915
$ws.Range("T4").Value = 0.0424684464222511
$ws.Range("W4").Value = 0.03486477573522739
$ws.Range("X4").Value = 0.004394526785374482
$ws.Range("Y4").Value = 0.003681724392147476
$ws.Range("Z4").Value = 0.04769709085191239
$ws.Range("AA4").Value = 0.01538347037523616
$ws.Range("AD4").Value = 0.0179976978033023
$ws.Range("AF4").Value = 0.01205174607016898
$ws.Range("AI4").Value = 0.01498004774979135
$ws.Range("D5").Value = 0.2818162333898037
$ws.Range("E5").Value = 0.1201554016982328
$ws.Range("F5").Value = 0.261023220444124
$ws.Range("H5").Value = 0.02075029777593906
$ws.Range("J5").Value = 0.008593943463216383
$ws.Range("K5").Value = 0.01089467913624654
$ws.Range("M5").Value = 0.01866724908379461
$ws.Range("N5").Value = 0.02698215942236624
$ws.Range("O5").Value = 0.08647845761436396
$ws.Range("S5").Value = 0.01517961559298407
$ws.Range("V5").Value = 0.005712868888609325
$ws.Range("X5").Value = 0.02458992941514733
$ws.Range("Z5").Value = 0.05505691056833358
$ws.Range("AA5").Value = 0.002399732318684702
$ws.Range("AB5").Value = 0.01356999317160426
$ws.Range("AD5").Value = 0.006842634444460041
$ws.Range("AE5").Value = 0.001983927539412212
$ws.Range("AF5").Value = 0.01747978806359695
$ws.Range("AI5").Value = 0.02182295796908004
$ws.Range("D6").Value = 0.2630237311841145
$ws.Range("E6").Value = 0.05386727994944034
$ws.Range("F6").Value = 0.3375679340521504
$ws.Range("H6").Value = 0.04076941166413879
$ws.Range("J6").Value = 0.004755764058608988
$ws.Range("M6").Value = 0.01923994433078173
$ws.Range("N6").Value = 0.01102200968272733
$ws.Range("O6").Value = 0.09964430582855854
$ws.Range("S6").Value = 0.009658899533697004
$ws.Range("T6").Value = 0.002155024117859209
$ws.Range("V6").Value = 0.00428031002709753
$ws.Range("X6").Value = 0.02934805327942772
$ws.Range("Z6").Value = 0.06529339219559213
$ws.Range("AA6").Value = 0.00363344543396361
$ws.Range("AB6").Value = 0.008704511804244201
$ws.Range("AD6").Value = 0.01231079178945854
$ws.Range("AF6").Value = 0.01343873810873514
$ws.Range("AI6").Value = 0.02128645295940439

# --- Step2_Sj (158 cell updates) ---
$ws = $wb.Worksheets.Item("Step2_Sj")
$ws.Range("D2").Value = 0.2286418605659838
$ws.Range("E2").Value = 0.2286418605659838
$ws.Range("F2").Value = 0.6762675483303523
$ws.Range("G2").Value = 0.6762675483303523
$ws.Range("H2").Value = 0.7022104516239559
$ws.Range("I2").Value = 0.7022104516239559
$ws.Range("J2").Value = 0.7022104516239559
$ws.Range("K2").Value = 0.7022104516239559
$ws.Range("L2").Value = 0.7022104516239559
$ws.Range("M2").Value = 0.717551302155228
$ws.Range("N2").Value = 0.7472984641446162
$ws.Range("O2").Value = 0.8268511476592505
$ws.Range("P2").Value = 0.8268511476592505
$ws.Range("Q2").Value = 0.8268511476592505
$ws.Range("R2").Value = 0.8268511476592505
$ws.Range("S2").Value = 0.873813915078804
$ws.Range("T2").Value = 0.8756583395272027
$ws.Range("U2").Value = 0.8756583395272027
$ws.Range("V2").Value = 0.8838039966015869
$ws.Range("W2").Value = 0.8838039966015869
$ws.Range("X2").Value = 0.9138538761796406
$ws.Range("Y2").Value = 0.9138538761796406
$ws.Range("Z2").Value = 0.9495674486796416
$ws.Range("AA2").Value = 0.9495674486796416
$ws.Range("AB2").Value = 0.9495674486796416
$ws.Range("AC2").Value = 0.9495674486796416
$ws.Range("AD2").Value = 0.9545222142585025
$ws.Range("AE2").Value = 0.9545222142585025
$ws.Range("AF2").Value = 0.9748623886123277
$ws.Range("AG2").Value = 0.9748623886123277
$ws.Range("AH2").Value = 0.9748623886123277
$ws.Range("AI2").Value = 1
$ws.Range("D3").Value = 0.2030359476101878
$ws.Range("E3").Value = 0.2030359476101878
$ws.Range("F3").Value = 0.670225842304407
$ws.Range("G3").Value = 0.6712237936768048
$ws.Range("H3").Value = 0.6973232753655244
$ws.Range("I3").Value = 0.6973232753655244
$ws.Range("J3").Value = 0.6973232753655244
$ws.Range("K3").Value = 0.6973232753655244
$ws.Range("L3").Value = 0.6973232753655244
$ws.Range("M3").Value = 0.7125403915177029
$ws.Range("N3").Value = 0.7125403915177029
$ws.Range("O3").Value = 0.8175328631617123
$ws.Range("P3").Value = 0.8175328631617123
$ws.Range("Q3").Value = 0.8175328631617123
$ws.Range("R3").Value = 0.8175328631617123
$ws.Range("S3").Value = 0.8256597462583475
$ws.Range("T3").Value = 0.8494228440498545
$ws.Range("U3").Value = 0.8494228440498545
$ws.Range("V3").Value = 0.8494228440498545
$ws.Range("W3").Value = 0.861595128675689
$ws.Range("X3").Value = 0.8731210284490267
$ws.Range("Y3").Value = 0.8731210284490267
$ws.Range("Z3").Value = 0.9184911645043582
$ws.Range("AA3").Value = 0.9338233638284976
$ws.Range("AB3").Value = 0.9338233638284976
$ws.Range("AC3").Value = 0.9338233638284976
$ws.Range("AD3").Value = 0.9505845754372002
$ws.Range("AE3").Value = 0.9505845754372002
$ws.Range("AF3").Value = 0.9747682271716367
$ws.Range("AG3").Value = 0.9747682271716367
$ws.Range("AH3").Value = 0.9747682271716367
$ws.Range("D4").Value = 0.1152529212304223
$ws.Range("E4").Value = 0.159386083900929
$ws.Range("F4").Value = 0.6218224851056837
$ws.Range("G4").Value = 0.6729763781689067
$ws.Range("H4").Value = 0.7122132441954077
$ws.Range("I4").Value = 0.7177107360081485
$ws.Range("J4").Value = 0.7177107360081485
$ws.Range("K4").Value = 0.7177107360081485
$ws.Range("L4").Value = 0.7177107360081485
$ws.Range("M4").Value = 0.7177107360081485
$ws.Range("N4").Value = 0.7177107360081485
$ws.Range("O4").Value = 0.8006293241612301
$ws.Range("P4").Value = 0.8006293241612301
$ws.Range("Q4").Value = 0.8006293241612301
$ws.Range("R4").Value = 0.8050793020501485
$ws.Range("S4").Value = 0.8064804738145884
$ws.Range("T4").Value = 0.8489489202368395
$ws.Range("U4").Value = 0.8489489202368395
$ws.Range("V4").Value = 0.8489489202368395
$ws.Range("W4").Value = 0.8838136959720669
$ws.Range("X4").Value = 0.8882082227574414
$ws.Range("Y4").Value = 0.8918899471495888
$ws.Range("Z4").Value = 0.9395870380015012
$ws.Range("AA4").Value = 0.9549705083767374
$ws.Range("AB4").Value = 0.9549705083767374
$ws.Range("AC4").Value = 0.9549705083767374
$ws.Range("AD4").Value = 0.9729682061800397
$ws.Range("AE4").Value = 0.9729682061800397
$ws.Range("AF4").Value = 0.9850199522502087
$ws.Range("AG4").Value = 0.9850199522502087
$ws.Range("AH4").Value = 0.9850199522502087
$ws.Range("AI4").Value = 1
$ws.Range("D5").Value = 0.2818162333898037
$ws.Range("E5").Value = 0.4019716350880365
$ws.Range("F5").Value = 0.6629948555321605
$ws.Range("G5").Value = 0.6629948555321605
$ws.Range("H5").Value = 0.6837451533080996
$ws.Range("I5").Value = 0.6837451533080996
$ws.Range("J5").Value = 0.692339096771316
$ws.Range("K5").Value = 0.7032337759075625
$ws.Range("L5").Value = 0.7032337759075625
$ws.Range("M5").Value = 0.7219010249913571
$ws.Range("N5").Value = 0.7488831844137234
$ws.Range("O5").Value = 0.8353616420280874
$ws.Range("P5").Value = 0.8353616420280874
$ws.Range("Q5").Value = 0.8353616420280874
$ws.Range("R5").Value = 0.8353616420280874
$ws.Range("S5").Value = 0.8505412576210715
$ws.Range("T5").Value = 0.8505412576210715
$ws.Range("U5").Value = 0.8505412576210715
$ws.Range("V5").Value = 0.8562541265096808
$ws.Range("W5").Value = 0.8562541265096808
$ws.Range("X5").Value = 0.8808440559248282
$ws.Range("Y5").Value = 0.8808440559248282
$ws.Range("Z5").Value = 0.9359009664931618
$ws.Range("AA5").Value = 0.9383006988118465
$ws.Range("AB5").Value = 0.9518706919834508
$ws.Range("AC5").Value = 0.9518706919834508
$ws.Range("AD5").Value = 0.9587133264279109
$ws.Range("AE5").Value = 0.9606972539673231
$ws.Range("AF5").Value = 0.97817704203092
$ws.Range("AG5").Value = 0.97817704203092
$ws.Range("AH5").Value = 0.97817704203092
$ws.Range("AI5").Value = 1
$ws.Range("D6").Value = 0.2630237311841145
$ws.Range("E6").Value = 0.3168910111335548
$ws.Range("F6").Value = 0.6544589451857052
$ws.Range("G6").Value = 0.6544589451857052
$ws.Range("H6").Value = 0.6952283568498441
$ws.Range("I6").Value = 0.6952283568498441
$ws.Range("J6").Value = 0.6999841209084531
$ws.Range("K6").Value = 0.6999841209084531
$ws.Range("L6").Value = 0.6999841209084531
$ws.Range("M6").Value = 0.7192240652392348
$ws.Range("N6").Value = 0.7302460749219621
$ws.Range("O6").Value = 0.8298903807505207
$ws.Range("P6").Value = 0.8298903807505207
$ws.Range("Q6").Value = 0.8298903807505207
$ws.Range("R6").Value = 0.8298903807505207
$ws.Range("S6").Value = 0.8395492802842177
$ws.Range("T6").Value = 0.8417043044020769
$ws.Range("U6").Value = 0.8417043044020769
$ws.Range("V6").Value = 0.8459846144291744
$ws.Range("W6").Value = 0.8459846144291744
$ws.Range("X6").Value = 0.8753326677086021
$ws.Range("Y6").Value = 0.8753326677086021
$ws.Range("Z6").Value = 0.9406260599041942
$ws.Range("AA6").Value = 0.9442595053381578
$ws.Range("AB6").Value = 0.9529640171424021
$ws.Range("AC6").Value = 0.9529640171424021
$ws.Range("AD6").Value = 0.9652748089318607
$ws.Range("AE6").Value = 0.9652748089318607
$ws.Range("AF6").Value = 0.9787135470405959
$ws.Range("AG6").Value = 0.9787135470405959
$ws.Range("AH6").Value = 0.9787135470405959

# --- Step3_DataPts_0.5 (6 cell updates) ---
$ws = $wb.Worksheets.Item("Step3_DataPts_0.5")
$ws.Range("F2").Value = 0.6762675483303523
$ws.Range("F3").Value = 0.670225842304407
$ws.Range("F4").Value = 0.6218224851056837
$ws.Range("D5").Value = 5
$ws.Range("G5").Value = 4
$ws.Range("F6").Value = 0.6544589451857052

# --- Step3_DataPts_0.7 (15 cell updates) ---
$ws = $wb.Worksheets.Item("Step3_DataPts_0.7")
$ws.Range("D2").Value = 7
$ws.Range("F2").Value = 0.7022104516239559
$ws.Range("G2").Value = 6
$ws.Range("D3").Value = 12
$ws.Range("F3").Value = 0.7125403915177029
$ws.Range("G3").Value = 11
$ws.Range("D4").Value = 7
$ws.Range("F4").Value = 0.7122132441954077
$ws.Range("G4").Value = 6
$ws.Range("D5").Value = 10
$ws.Range("F5").Value = 0.7032337759075625
$ws.Range("G5").Value = 9
$ws.Range("D6").Value = 12
$ws.Range("F6").Value = 0.7192240652392348
$ws.Range("G6").Value = 11

# --- Step3_DataPts_0.8 (15 cell updates) ---
$ws = $wb.Worksheets.Item("Step3_DataPts_0.8")
$ws.Range("D2").Value = 14
$ws.Range("F2").Value = 0.8268511476592505
$ws.Range("G2").Value = 13
$ws.Range("D3").Value = 14
$ws.Range("F3").Value = 0.8175328631617123
$ws.Range("G3").Value = 13
$ws.Range("D4").Value = 14
$ws.Range("F4").Value = 0.8006293241612301
$ws.Range("G4").Value = 13
$ws.Range("D5").Value = 14
$ws.Range("F5").Value = 0.8353616420280874
$ws.Range("G5").Value = 13
$ws.Range("D6").Value = 14
$ws.Range("F6").Value = 0.8298903807505207
$ws.Range("G6").Value = 13

# --- Step3_DataPts_0.9 (15 cell updates) ---
$ws = $wb.Worksheets.Item("Step3_DataPts_0.9")
$ws.Range("D2").Value = 23
$ws.Range("F2").Value = 0.9138538761796406
$ws.Range("G2").Value = 22
$ws.Range("D3").Value = 25
$ws.Range("F3").Value = 0.9184911645043582
$ws.Range("G3").Value = 24
$ws.Range("D4").Value = 25
$ws.Range("F4").Value = 0.9395870380015012
$ws.Range("G4").Value = 24
$ws.Range("D5").Value = 25
$ws.Range("F5").Value = 0.9359009664931618
$ws.Range("G5").Value = 24
$ws.Range("D6").Value = 25
$ws.Range("F6").Value = 0.9406260599041942
$ws.Range("G6").Value = 24
